$d = $word.ActiveDocument

# Locate the "Version 1." text so we work off real character offsets
# rather than hard-coded ones.
$found = $d.Content.Duplicate
[void]$found.Find.Execute("Version 1.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $found.Start

# 1. Split the single "Version" run into "Versi" + "on" runs (no run
#    formatting is introduced) by re-inserting the "on" substring as its
#    own run via InsertXML. Word then splits the surrounding run around
#    it while leaving the proofErr/bookmark markers in place.
$onRng = $d.Range($base + 5, $base + 7)
$onRng.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>on</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 2. Change " 1." to " 2" (the trailing period is re-added separately below,
#    after the bookmark, as its own run).
$verRng = $d.Range($base + 8, $base + 10)
$verRng.Text = "2"

# 3. Append a new "." run at the very end of the paragraph (i.e. after the
#    _GoBack bookmark).
$paraEnd = $d.Paragraphs(1).Range.End
$endRng = $d.Range($paraEnd, $paraEnd)
$endRng.InsertAfter(".")
